# Jun's file updates for all IO data and others
#
# 1. "OECD Data" sheet: the "National GDP" column (G) header is renamed to
#    "Value", and the helper "National Growth" (H) / "State GDP" (I) columns
#    (headers + all per-year formulas) are removed entirely.
# 2. "BGDP" sheet: the yearly GDP formulas in column B no longer pull the
#    (now-removed) State GDP column I from "OECD Data"; instead they compute
#    the value directly from the National GDP column G, converted using the
#    factors kept on the "About" sheet (A15 = 2012-dollars-per-2010-dollar,
#    A16 = dollars-per-million-dollars).
# 3. The workbook re-opens with the "About" sheet active/selected instead of
#    "BGDP", and the stale K6/E43 selections on "OECD Data"/"BGDP" reset to A1.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsOECD  = $wb.Worksheets.Item("OECD Data")
$wsBGDP  = $wb.Worksheets.Item("BGDP")

# --- "OECD Data" sheet: collapse National Growth / State GDP helper columns ---

# Rename the remaining "National GDP" header to "Value".
$wsOECD.Range("G1").Value = "Value"

# Remove the "National Growth" (H) and "State GDP" (I) columns: header cells
# in row 1 plus every formula in rows 2-48.
$wsOECD.Range("H1:I48").ClearContents()

# --- "BGDP" sheet: recompute GDP straight from National GDP (G) on "OECD Data" ---

for ($r = 2; $r -le 48; $r++) {
    $wsBGDP.Cells.Item($r, 2).Formula = '=''OECD Data''!G' + $r + '*About!$A$15*About!$A$16'
}

# --- Selection / active-sheet bookkeeping (matches the saved view state) ---

$wsOECD.Activate() | Out-Null
$wsOECD.Range("A1").Select() | Out-Null

$wsBGDP.Activate() | Out-Null
$wsBGDP.Range("A1").Select() | Out-Null

$wsAbout.Activate() | Out-Null
$wsAbout.Range("A1").Select() | Out-Null
